# Add 2022-Q3 data:
#  - insert a new worksheet "2022-Q3" right after "总计", before "2022-Q2"
#  - populate it with the new quarter's fund-holding table
#  - prepend the new quarter's summary row into "总计", shifting the
#    existing quarterly summary rows down by one

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)      # 总计
$q2Sheet    = $wb.Worksheets.Item(2)      # 2022-Q2 (template for the new sheet's layout)

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by copying "2022-Q2" (same columns,
#    header texts, styles, borders) right after "总计", then overwrite
#    its data with the 2022-Q3 numbers.
# ---------------------------------------------------------------------
$q2Sheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Fund code column (B) and the text-valued numeric columns (D:G) need to
# stay text (leading zeros such as 013273 / 015674 must not turn numeric).
$q3Sheet.Range("B2:B9").NumberFormat = "@"
$q3Sheet.Range("D2:G9").NumberFormat = "@"

$q3Data = @(
    @("013273", "招商沪深300地产等权重指数C",   "6.78",  "94.58", "10.94", "0.7417", 6),
    @("512200", "南方中证全指房地产ETF",         "40.02", "99.99", "1.85",  "0.7404", 9),
    @("161721", "招商沪深300地产等权重指数A",    "4.91",  "94.58", "10.94", "0.5372", 6),
    @("515060", "华夏中证全指房地产ETF",         "4.09",  "99.07", "1.83",  "0.0748", 10),
    @("160628", "鹏华中证800地产指数（LOF）A",   "2.70",  "94.64", "2.76",  "0.0745", 9),
    @("159768", "银华中证内地地产主题ETF",       "1.78",  "97.47", "2.86",  "0.0509", 9),
    @("159707", "华宝中证800地产ETF",            "1.17",  "98.55", "2.86",  "0.0335", 9),
    @("015674", "鹏华中证800地产指数（LOF）C",   "0.91",  "94.64", "2.76",  "0.0251", 9)
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $row = $i + 2
    $item = $q3Data[$i]
    $q3Sheet.Cells.Item($row, 2).Value = $item[0]
    $q3Sheet.Cells.Item($row, 3).Value = $item[1]
    $q3Sheet.Cells.Item($row, 4).Value = $item[2]
    $q3Sheet.Cells.Item($row, 5).Value = $item[3]
    $q3Sheet.Cells.Item($row, 6).Value = $item[4]
    $q3Sheet.Cells.Item($row, 7).Value = $item[5]
    $q3Sheet.Cells.Item($row, 8).Value = $item[6]
}

# ---------------------------------------------------------------------
# 2) "总计" sheet: add the new row9 (carrying the old row8's data down),
#    formatted like the existing A-column cells, then shift the B:D
#    content of every row down by one and place the new 2022-Q3 summary
#    values into row 2.
# ---------------------------------------------------------------------
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$totalData = @(
    @("2022-Q3", 8,  2.28),
    @("2022-Q2", 8,  4.38),
    @("2022-Q1", 22, 6.73),
    @("2021-Q4", 7,  4.99),
    @("2021-Q3", 10, 4.48),
    @("2021-Q2", 6,  5.07),
    @("2021-Q1", 9,  5.57),
    @("2020-Q4", 10, 6.03)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $row = $i + 2
    $item = $totalData[$i]
    $totalSheet.Cells.Item($row, 1).Value = $i
    $totalSheet.Cells.Item($row, 2).Value = $item[0]
    $totalSheet.Cells.Item($row, 3).Value = $item[1]
    $totalSheet.Cells.Item($row, 4).Value = $item[2]
}
